# Add "Master Cabang" block to the "Coding Progress" sheet, mirroring the
# existing "Master Group" / "Master User" blocks (rows 9-11 below the
# existing data that ends at row 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coding Progress")

$ws.Range("B9").Value = "Master Cabang"
$ws.Range("C9").Value = "DONE"
$ws.Range("D9").Value = "Data Entry Coding"

$ws.Range("C10").Value = "partial DONE"
$ws.Range("D10").Value = "Input validation"

$ws.Range("C11").Value = "DONE"
$ws.Range("D11").Value = "DataGrid for browsing"

# Match the author's resulting selection (rows 9-11, columns C:D) and keep
# the active cell at C9, as seen in the saved workbook.
$ws.Range("C9:D11").Select()
